# Generate Report for Handback
# Update the handback-status report timestamps / priority value that were
# regenerated by the CI job.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for 57779544-... (Overview, rows 3 & 5
# shared the same value and still do after the refresh).
$wsOverview.Range("G3").Value = "2016-08-31 08:18:53"
$wsOverview.Range("G5").Value = "2016-08-31 08:18:53"

# "Priority" value - shared by zh-cn and de-de sheets, rows for
# 57779544-... (row 3) and d90ab85e-... (row 5).
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# zh-cn "Correspond Handoff Datetime"
$wsZhCn.Range("H3").Value = "2016-08-31 08:18:47"
$wsZhCn.Range("H5").Value = "2016-08-31 08:18:47"

# zh-cn "Correspond Handback DateTime"
$wsZhCn.Range("K3").Value = "2016-08-31 08:19:22"
$wsZhCn.Range("K5").Value = "2016-08-31 08:19:22"

# de-de "Correspond Handoff Datetime" (shares the value with Overview G3/G5)
$wsDeDe.Range("H3").Value = "2016-08-31 08:18:53"
$wsDeDe.Range("H5").Value = "2016-08-31 08:18:53"

# de-de "Correspond Handback DateTime"
$wsDeDe.Range("K3").Value = "2016-08-31 08:19:30"
$wsDeDe.Range("K5").Value = "2016-08-31 08:19:30"
